# Introduction.pptx update
#  1. Slide 9 title: "Introduction" -> "Features"
#  2. Refresh the auto "datetimeFigureOut" date placeholders (slide master,
#     every slide layout and the notes master) from 10/21/2020 -> 10/25/2020,
#     which is what PowerPoint does to "Update automatically" date fields
#     whenever the deck is saved on a later day.

$p = $ppt.ActivePresentation

# --- 1. Slide title -------------------------------------------------------
$slide = $p.Slides.Item(9)
$titleShape = $slide.Shapes.Item(1)
if ($titleShape.HasTextFrame) {
    if ($titleShape.TextFrame.TextRange.Text -eq "Introduction") {
        $titleShape.TextFrame.TextRange.Text = "Features"
    }
}

# --- 2. Date placeholders ---------------------------------------------------
$oldDate = "10/21/2020"
$newDate = "10/25/2020"

function Update-DatePlaceholder {
    param($shapes)

    foreach ($sh in $shapes) {
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout that hangs off the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes
